$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 79
$ws.Range("E2").Value = 80
$ws.Range("F2").Value = 2420
$ws.Range("G2").Value = 815
$ws.Range("H2").Value = 811
$ws.Range("I2").Value = 387
$ws.Range("J2").Value = 528
$ws.Range("K2").Value = 15281
$ws.Range("L2").Value = 548
$ws.Range("M2").Value = 495
$ws.Range("N2").Value = 788
$ws.Range("O2").Value = 6.936708860759493
$ws.Range("P2").Value = 6.265822784810126
$ws.Range("Q2").Value = 9.974683544303797
$ws.Range("R2").Value = 10.32
$ws.Range("S2").Value = 1837.97

$ws.Range("D3").Value = 41
$ws.Range("E3").Value = 22
$ws.Range("F3").Value = 1171
$ws.Range("G3").Value = 442
$ws.Range("H3").Value = 450
$ws.Range("I3").Value = 84
$ws.Range("J3").Value = 93
$ws.Range("K3").Value = 4954
$ws.Range("L3").Value = 112
$ws.Range("M3").Value = 187
$ws.Range("N3").Value = 300
$ws.Range("O3").Value = 2.731707317073171
$ws.Range("P3").Value = 4.560975609756097
$ws.Range("Q3").Value = 7.317073170731708
$ws.Range("R3").Value = 10.78
$ws.Range("S3").Value = 1713.66

$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 296
$ws.Range("E4").Value = 171
$ws.Range("F4").Value = 8648
$ws.Range("G4").Value = 2552
$ws.Range("H4").Value = 2563
$ws.Range("I4").Value = 868
$ws.Range("J4").Value = 795
$ws.Range("K4").Value = 58854
$ws.Range("L4").Value = 2091
$ws.Range("M4").Value = 1638
$ws.Range("N4").Value = 2261
$ws.Range("O4").Value = 7.064189189189189
$ws.Range("P4").Value = 5.533783783783784
$ws.Range("Q4").Value = 7.638513513513513
$ws.Range("R4").Value = 8.619999999999999
$ws.Range("S4").Value = 1752.97

$ws.Range("D5").Value = 131
$ws.Range("E5").Value = 104
$ws.Range("F5").Value = 4022
$ws.Range("G5").Value = 2575
$ws.Range("H5").Value = 2508
$ws.Range("I5").Value = 675
$ws.Range("J5").Value = 732
$ws.Range("K5").Value = 21160
$ws.Range("L5").Value = 642
$ws.Range("M5").Value = 449
$ws.Range("N5").Value = 1311
$ws.Range("O5").Value = 4.900763358778626
$ws.Range("P5").Value = 3.427480916030534
$ws.Range("Q5").Value = 10.00763358778626
$ws.Range("R5").Value = 19.66
$ws.Range("S5").Value = 1842.14

$ws.Range("F6").Value = 1421
$ws.Range("G6").Value = 1044
$ws.Range("J6").Value = 354
$ws.Range("R6").Value = 21.31
$ws.Range("S6").Value = 1740

$ws.Range("D7").Value = 21
$ws.Range("E7").Value = 17
$ws.Range("F7").Value = 626
$ws.Range("G7").Value = 241
$ws.Range("H7").Value = 254
$ws.Range("I7").Value = 45
$ws.Range("J7").Value = 57
$ws.Range("K7").Value = 4532
$ws.Range("L7").Value = 83
$ws.Range("M7").Value = 79
$ws.Range("N7").Value = 120
$ws.Range("O7").Value = 3.952380952380953
$ws.Range("P7").Value = 3.761904761904762
$ws.Range("Q7").Value = 5.714285714285714
$ws.Range("R7").Value = 11.48
$ws.Range("S7").Value = 1788.57

$ws.Range("D8").Value = 57
$ws.Range("E8").Value = 22
$ws.Range("F8").Value = 1652
$ws.Range("G8").Value = 2191
$ws.Range("H8").Value = 2053
$ws.Range("I8").Value = 683
$ws.Range("J8").Value = 439
$ws.Range("K8").Value = 2069
$ws.Range("L8").Value = 157
$ws.Range("M8").Value = 411
$ws.Range("N8").Value = 819
$ws.Range("O8").Value = 2.754385964912281
$ws.Range("P8").Value = 7.210526315789473
$ws.Range("Q8").Value = 14.36842105263158
$ws.Range("R8").Value = 38.44
$ws.Range("S8").Value = 1738.95

$ws.Range("F9").Value = 1634
$ws.Range("G9").Value = 641
$ws.Range("J9").Value = 158
$ws.Range("R9").Value = 11.87
$ws.Range("S9").Value = 1815.56

$ws.Range("D10").Value = 121
$ws.Range("E10").Value = 70
$ws.Range("F10").Value = 3637
$ws.Range("G10").Value = 4159
$ws.Range("H10").Value = 4178
$ws.Range("I10").Value = 986
$ws.Range("J10").Value = 1002
$ws.Range("K10").Value = 9461
$ws.Range("L10").Value = 483
$ws.Range("M10").Value = 391
$ws.Range("N10").Value = 1661
$ws.Range("O10").Value = 3.991735537190082
$ws.Range("P10").Value = 3.231404958677686
$ws.Range("Q10").Value = 13.72727272727273
$ws.Range("R10").Value = 34.37
$ws.Range("S10").Value = 1803.47

$ws.Range("D11").Value = 49
$ws.Range("E11").Value = 32
$ws.Range("F11").Value = 1517
$ws.Range("G11").Value = 685
$ws.Range("H11").Value = 670
$ws.Range("I11").Value = 141
$ws.Range("J11").Value = 251
$ws.Range("K11").Value = 9411
$ws.Range("L11").Value = 309
$ws.Range("M11").Value = 250
$ws.Range("N11").Value = 381
$ws.Range("O11").Value = 6.306122448979592
$ws.Range("P11").Value = 5.102040816326531
$ws.Range("Q11").Value = 7.775510204081633
$ws.Range("R11").Value = 13.98
$ws.Range("S11").Value = 1857.55

$ws.Range("F12").Value = 400
$ws.Range("G12").Value = 116
$ws.Range("J12").Value = 4
$ws.Range("R12").Value = 8.92
$ws.Range("S12").Value = 1846.15

$ws.Range("D13").Value = 27
$ws.Range("E13").Value = 13
$ws.Range("F13").Value = 769
$ws.Range("G13").Value = 317
$ws.Range("H13").Value = 336
$ws.Range("I13").Value = 86
$ws.Range("J13").Value = 107
$ws.Range("K13").Value = 5394
$ws.Range("L13").Value = 112
$ws.Range("M13").Value = 86
$ws.Range("N13").Value = 263
$ws.Range("O13").Value = 4.148148148148148
$ws.Range("P13").Value = 3.185185185185185
$ws.Range("Q13").Value = 9.74074074074074
$ws.Range("R13").Value = 11.74
$ws.Range("S13").Value = 1708.89
